$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56 (existing rows 56-65 shift down to 57-66)
$ws.Rows("56:56").Insert()

# Populate the newly inserted row 56 with the new weekly price-report entry
$ws.Range("A56").Value = 9
$ws.Range("B56").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C56").Value = "Metropolitana"
$ws.Range("D56").Value = 44694
$ws.Range("E56").Value = 13
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100102
$ws.Range("H56").Value = "Cítricos"
$ws.Range("I56").Value = 100102006
$ws.Range("J56").Value = "Pomelo"
$ws.Range("K56").Value = "Start Ruby"
$ws.Range("L56").Value = "Primera"
$ws.Range("M56").Value = 350
$ws.Range("N56").Value = 7500
$ws.Range("O56").Value = 7500
$ws.Range("P56").Value = 7500
$ws.Range("Q56").Value = "$/caja 14 kilos granel"
$ws.Range("R56").Value = "Región Metropolitana"
$ws.Range("S56").Value = 536
$ws.Range("T56").Value = 14
